$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data for rows 2-15: Word, X, Y, Width, Height
$data = @(
    @("19th ", 492, 644, 42, 23),
    @("century. ", 534, 644, 75, 23),
    @([string]([char]0x201C) + "Romantic ", 757, 644, 96, 23),
    @("period" + [string]([char]0x201D) + " ", 853, 644, 68, 23),
    @("the ", 669, 669, 32, 23),
    @("Classical ", 701, 669, 79, 23),
    @("period, ", 780, 669, 65, 23),
    @("18th-century ", 287, 735.4, 114, 23),
    @("Pastoral," + [string]([char]0x201D) + " ", 454, 1034.6, 88, 23),
    @("sea ", 268, 1126, 32, 23),
    @("coming ", 300, 1126, 68, 23),
    @("into ", 368, 1126, 39, 23),
    @("Fingal" + [string]([char]0x2019) + "s ", 407, 1126, 71, 23),
    @("Cave ", 478, 1126, 46, 23)
)

$row = 2
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $ws.Cells.Item($row, 4).Value = $entry[3]
    $ws.Cells.Item($row, 5).Value = $entry[4]
    $row++
}
